$d = $word.ActiveDocument

# 1. Footer page-number field's cached result: "3" -> "8"
$ftrRange = $d.Sections(1).Footers(1).Range
$ftrRange.Find.Execute("3", $true, $true, $false, $false, $false, $true, 1, $false, "8", 2)

# 2. Rename character style "LineNumbering": "Line Numbering" -> "Line Number"
$d.Styles("LineNumbering").NameLocal = "Line Number"

# 3. Rename character style "FootnoteAnchor": "Footnote Anchor" -> "Footnote Reference"
$d.Styles("FootnoteAnchor").NameLocal = "Footnote Reference"

# 4. Rename character style "EndnoteAnchor": "Endnote Anchor" -> "Endnote Reference"
$d.Styles("EndnoteAnchor").NameLocal = "Endnote Reference"

Write-Output "edits applied"
